# ONSDeathStats.xlsx edit: add "England and Wales" analysis sheet
# 1) Rename existing Sheet1 -> "Nottinghamshire"
# 2) Insert a new sheet "England and Wales" after it, populate with
#    monthly-deaths data for 2015-2020, summary stats and make it the
#    active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1. rename the original sheet -----------------------------------------
$wsNotts = $wb.Worksheets.Item(1)
$wsNotts.Name = "Nottinghamshire"

# --- 2. add the new sheet after it -----------------------------------------
$ws = $wb.Worksheets.Add($null, $wsNotts)
$ws.Name = "England and Wales"

$url = "https://www.ons.gov.uk/peoplepopulationandcommunity/birthsdeathsandmarriages/deaths/datasets/monthlyfiguresondeathsregisteredbyareaofusualresidence"

# --- column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.33
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Range("C1:K1").EntireColumn.ColumnWidth = 6.55
$ws.Columns.Item(16).ColumnWidth = 12

# --- row 1: source link ------------------------------------------------------
$ws.Range("A1").Value = $url
$ws.Hyperlinks.Add($ws.Range("A1"), $url) | Out-Null

# --- row 3: month headers ----------------------------------------------------
$months = New-Object 'object[,]' 1,12
$months[0,0]  = "Jan"
$months[0,1]  = "Feb"
$months[0,2]  = "Mar"
$months[0,3]  = "Apr"
$months[0,4]  = "May"
$months[0,5]  = "Jun"
$months[0,6]  = "Jul"
$months[0,7]  = "Aug"
$months[0,8]  = "Sep"
$months[0,9]  = "Oct"
$months[0,10] = "Nov"
$months[0,11] = "Dec"
$ws.Range("C3:N3").Value = $months
$ws.Range("P3").Value = "Total to Sep"
$ws.Range("C3:N3").Font.Bold = $true
$ws.Range("P3").Font.Bold = $true

# --- rows 4-9: yearly data ----------------------------------------------------
$years = @(
  @{ Row=4; Label="ENGLAND AND WALES 2020"; Data=@(56597,43555,49641,88049,52315,42577,40731,37129,42432) },
  @{ Row=5; Label="ENGLAND AND WALES 2019"; Data=@(53774,45695,43817,44005,44292,38511,42192,38721,39915,46131,45124,47376) },
  @{ Row=6; Label="ENGLAND AND WALES 2018"; Data=@(64020,49087,51131,46383,42685,39679,40621,40071,37013,44311,43834,41430) },
  @{ Row=7; Label="ENGLAND AND WALES 2017"; Data=@(57266,47695,48577,39024,44183,42074,38314,40963,40002,43504,45476,45052) },
  @{ Row=8; Label="ENGLAND AND WALES 2016"; Data=@(47351,45922,48562,46755,41291,41921,38882,40676,40250,40360,46418,45469) },
  @{ Row=9; Label="ENGLAND AND WALES 2015"; Data=@(60779,46634,47820,45077,39250,41992,40400,36096,41491,42125,41431,45412) }
)

foreach ($y in $years) {
    $r = $y.Row
    $data = $y.Data
    $ws.Range("A$r").Value = "K04000001"
    $ws.Range("B$r").Value = $y.Label
    $ws.Range("A$r" + ":B$r").Font.Bold = $true

    $n = $data.Length
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) { $arr[0,$i] = $data[$i] }
    $startCol = 3
    $endCol = $startCol + $n - 1
    $startA1 = $ws.Cells.Item($r, $startCol).Address($false, $false)
    $endA1   = $ws.Cells.Item($r, $endCol).Address($false, $false)
    $ws.Range("$startA1" + ":" + "$endA1").Value = $arr

    $ws.Range("P$r").Formula = "=SUM(C$r" + ":K$r)"
}

# highlight the Apr/May 2020 spike (columns F & G) in red, as in the source data
$ws.Range("F4:G4").Font.Color = 255
$ws.Range("P4").Font.Color = 255

# --- row 10: blank separator (kept un-styled) --------------------------------

# --- row 11: 5 year mean excluding 2020 --------------------------------------
$ws.Range("B11").Value = "5 Year Mean Excluding 2020"
$ws.Range("B11").Font.Bold = $true
$ws.Range("B11").HorizontalAlignment = -4131
$ws.Range("C11:N11").Formula = "=AVERAGE(C5:C9)"
$ws.Range("P11").Formula = "=AVERAGE(P5:P9)"

# --- row 12: 5 year std dev excluding 2020 -----------------------------------
$ws.Range("B12").Value = "5 Year Std Dev Excluding 2020"
$ws.Range("B12").Font.Bold = $true
$ws.Range("B12").HorizontalAlignment = -4131
$ws.Range("C12:N12").Formula = "=STDEV.S(C5:C9)"
$ws.Range("P12").Formula = "=STDEV.S(P5:P9)"

# --- row 13: blank ------------------------------------------------------------

# --- row 14: 2020 excess over 5 year mean ------------------------------------
$ws.Range("B14").Value = "2020 Excess over 5 Year Mean"
$ws.Range("C14:K14").Formula = "=C4-C11"
$ws.Range("C14:K14").NumberFormat = "#,##0"
$ws.Range("P14").Formula = "=P4-P11"
$ws.Range("P14").NumberFormat = "#,##0"
$ws.Range("F14:G14").Font.Color = 255
$ws.Range("P14").Font.Color = 255

# --- row 15: number of std deviations ----------------------------------------
$ws.Range("B15").Value = "Number of Std Deviations"
$ws.Range("C15:K15").Formula = "=C14/C12"
$ws.Range("P15").Formula = "=P14/P12"
$ws.Range("F15:G15").Font.Color = 255
$ws.Range("P15").Font.Color = 255

# --- activate the new sheet/selection ----------------------------------------
$wsNotts.Range("A1").Select() | Out-Null
$ws.Activate()
$ws.Range("O21").Select() | Out-Null

Write-Host "done"
